# Update "想去人数" (interested-count) figures in column F for the
# "展览" and "全部类型" sheets, matching the latest scrape snapshot.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# 展览 (sheet 1): row -> new F value
$wsExhibit.Range("F6").Value  = 30
$wsExhibit.Range("F7").Value  = 580
$wsExhibit.Range("F8").Value  = 65
$wsExhibit.Range("F9").Value  = 8421
$wsExhibit.Range("F10").Value = 786
$wsExhibit.Range("F11").Value = 310
$wsExhibit.Range("F12").Value = 1128
$wsExhibit.Range("F13").Value = 896
$wsExhibit.Range("F14").Value = 74
$wsExhibit.Range("F15").Value = 42
$wsExhibit.Range("F17").Value = 167
$wsExhibit.Range("F19").Value = 223
$wsExhibit.Range("F20").Value = 942

# 全部类型 (sheet 4): row -> new F value
$wsAll.Range("F7").Value  = 30
$wsAll.Range("F9").Value  = 580
$wsAll.Range("F10").Value = 65
$wsAll.Range("F11").Value = 8421
$wsAll.Range("F12").Value = 786
$wsAll.Range("F13").Value = 310
$wsAll.Range("F14").Value = 1128
$wsAll.Range("F15").Value = 896
$wsAll.Range("F16").Value = 74
$wsAll.Range("F17").Value = 42
$wsAll.Range("F19").Value = 167
$wsAll.Range("F21").Value = 223
$wsAll.Range("F22").Value = 942
